# Actualizacion automatica 2025-10-06 14:53:36
#
# Inserts two new client rows ("CERAMICAS Y MATERIALES LA ECONOMIA
# PENINSULAR S.A. CERMAPENSA" and "JACOME MONCAYO JAVIER ALFONSO") into the
# alphabetically sorted client tables on the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, widens column B, and refreshes the trailing
# summary row on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R, currency-style data C:R)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert "CERAMICAS ..." before "COMERCIAL INTERNACIONAL ..." (old row 14)
$ws1.Rows.Item(14).Insert()
$ws1.Cells.Item(14, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws1.Cells.Item(14, 2).Value = "CERAMICAS Y MATERIALES LA ECONOMIA PENINSULAR S.A. CERMAPENSA"
$ws1.Cells.Item(14, 3).Value = 259.2
$ws1.Cells.Item(14, 4).Value = 0
$ws1.Cells.Item(14, 5).Value = 0
$ws1.Cells.Item(14, 6).Value = 0
$ws1.Cells.Item(14, 7).Value = 0
$ws1.Cells.Item(14, 8).Value = 0
$ws1.Cells.Item(14, 9).Value = 0
$ws1.Cells.Item(14, 10).Value = 0
$ws1.Cells.Item(14, 11).Value = 0
$ws1.Cells.Item(14, 12).Value = 648.83
$ws1.Cells.Item(14, 13).Value = -43.78
$ws1.Cells.Item(14, 14).Value = 0
$ws1.Cells.Item(14, 15).Value = 0
$ws1.Cells.Item(14, 16).Value = 0
$ws1.Cells.Item(14, 17).Value = 0
$ws1.Cells.Item(14, 18).Value = 0

# Insert "JACOME MONCAYO JAVIER ALFONSO" before "JOWIN SA" (old row 30,
# now row 31 after the first insert above)
$ws1.Rows.Item(31).Insert()
$ws1.Cells.Item(31, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws1.Cells.Item(31, 2).Value = "JACOME MONCAYO JAVIER ALFONSO"
$ws1.Cells.Item(31, 3).Value = 0
$ws1.Cells.Item(31, 4).Value = 0
$ws1.Cells.Item(31, 5).Value = 0
$ws1.Cells.Item(31, 6).Value = 0
$ws1.Cells.Item(31, 7).Value = 0
$ws1.Cells.Item(31, 8).Value = 0
$ws1.Cells.Item(31, 9).Value = 0
$ws1.Cells.Item(31, 10).Value = 0
$ws1.Cells.Item(31, 11).Value = 0
$ws1.Cells.Item(31, 12).Value = 0
$ws1.Cells.Item(31, 13).Value = 0
$ws1.Cells.Item(31, 14).Value = 0
$ws1.Cells.Item(31, 15).Value = 0
$ws1.Cells.Item(31, 16).Value = 0
$ws1.Cells.Item(31, 17).Value = 0
$ws1.Cells.Item(31, 18).Value = 0

# Widen client-name column
$ws1.Columns.Item(2).ColumnWidth = 62.1667

# Refresh the "X de 57" summary row, now shifted down to row 59
$ws1.Cells.Item(59, 3).Value = "1 de 57"
$ws1.Cells.Item(59, 4).Value = "0 de 57"
$ws1.Cells.Item(59, 5).Value = "0 de 57"
$ws1.Cells.Item(59, 6).Value = "0 de 57"
$ws1.Cells.Item(59, 7).Value = "0 de 57"
$ws1.Cells.Item(59, 8).Value = "1 de 57"
$ws1.Cells.Item(59, 9).Value = "1 de 57"
$ws1.Cells.Item(59, 10).Value = "0 de 57"
$ws1.Cells.Item(59, 11).Value = "0 de 57"
$ws1.Cells.Item(59, 12).Value = "1 de 57"
$ws1.Cells.Item(59, 13).Value = "0 de 57"
$ws1.Cells.Item(59, 14).Value = "0 de 57"
$ws1.Cells.Item(59, 15).Value = "0 de 57"
$ws1.Cells.Item(59, 16).Value = "1 de 57"
$ws1.Cells.Item(59, 17).Value = "0 de 57"
$ws1.Cells.Item(59, 18).Value = "0 de 57"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G, monthly sales + PRESUPUESTO)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert "CERAMICAS ..." before "COMERCIAL INTERNACIONAL ..." (old row 14)
$ws2.Rows.Item(14).Insert()
$ws2.Cells.Item(14, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws2.Cells.Item(14, 2).Value = "CERAMICAS Y MATERIALES LA ECONOMIA PENINSULAR S.A. CERMAPENSA"
$ws2.Cells.Item(14, 3).Value = 4564.61
$ws2.Cells.Item(14, 4).Value = 11261.15
$ws2.Cells.Item(14, 5).Value = 8223.31
$ws2.Cells.Item(14, 6).Value = 864.25
$ws2.Cells.Item(14, 7).Value = 0

# Insert "JACOME MONCAYO JAVIER ALFONSO" before "JOWIN SA" (old row 30,
# now row 31 after the first insert above)
$ws2.Rows.Item(31).Insert()
$ws2.Cells.Item(31, 1).Value = "CASTRO ALCIVAR EDA MARIA"
$ws2.Cells.Item(31, 2).Value = "JACOME MONCAYO JAVIER ALFONSO"
$ws2.Cells.Item(31, 3).Value = 6665.35
$ws2.Cells.Item(31, 4).Value = 11275.94
$ws2.Cells.Item(31, 5).Value = 8872.17
$ws2.Cells.Item(31, 6).Value = 0
$ws2.Cells.Item(31, 7).Value = 0

# Widen client-name column
$ws2.Columns.Item(2).ColumnWidth = 62.1667

# Refresh the totals row, now shifted down to row 59
$ws2.Cells.Item(59, 3).Value = 94935.39
$ws2.Cells.Item(59, 4).Value = 81440.64
$ws2.Cells.Item(59, 5).Value = 91039.07
$ws2.Cells.Item(59, 6).Value = 1535.18
$ws2.Cells.Item(59, 7).Value = 85274.88
